# Generate Report for Handoff
# Replaces the old GUID-named file references (c68f88e4-...) with the new
# ones (a5c131ed-...) across all three sheets, updates the related
# timestamps, and refreshes the hyperlink display text while keeping the
# original hyperlink target (URL) intact.

$wb = $excel.ActiveWorkbook

$oldGuid = "c68f88e4-19cc-4d5e-8358-19d9c3c16285"
$newGuid = "a5c131ed-efda-418f-80f6-33a7da52d47e"

$oldHash = "b6615fb97375061884b39485e7e8a65409141eff"
$newHash = "fe6262f96bb8ccd23f1885eec2bb6c73eb816dfe"

$newFileName       = "$newGuid.md"
$newPathAndName    = "e2e\$newGuid.md"
$newHoDate         = "2016-09-01 09:14:36"

$newZhHandoffFile  = "$newGuid.$newHash.zh-cn.xlf"
$newZhHandoffDate  = "2016-09-01 09:14:31"

$newDeHandoffFile  = "$newGuid.$newHash.de-de.xlf"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb41676931de77368ac6f8eba6557ff39da7008b/e2e/$oldGuid.md"

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, $newPathAndName)

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhHandoffFile
$wsZh.Range("H2").Value = $newZhHandoffDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, $newFileName)

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeHandoffFile
$wsDe.Range("H2").Value = $newHoDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkAddress, [System.Type]::Missing, [System.Type]::Missing, $newFileName)
